$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.2138157894736842
$ws.Range("C2").Value2 = 0.4835526315789473
$ws.Range("J2").Value2 = 0.02960526315789474
$ws.Range("P2").Value2 = 0.1447368421052632
$ws.Range("S2").Value2 = 0.1282894736842105
$ws.Range("B3").Value2 = 0.02
$ws.Range("C3").Value2 = 0.02
$ws.Range("J3").Value2 = 0.04
$ws.Range("P3").Value2 = 0.7333333333333333
$ws.Range("S3").Value2 = 0.1866666666666667
$ws.Range("J4").Value2 = 0.07692307692307693
$ws.Range("P4").Value2 = 0.5769230769230769
$ws.Range("S4").Value2 = 0.3461538461538461
$ws.Range("P5").Value2 = 1
$ws.Range("B6").Value2 = 0.0625
$ws.Range("D6").Value2 = 0.01785714285714286
$ws.Range("E6").Value2 = 0.004464285714285714
$ws.Range("F6").Value2 = 0.04910714285714286
$ws.Range("J6").Value2 = 0.2008928571428572
$ws.Range("O6").Value2 = 0.01339285714285714
$ws.Range("Q6").Value2 = 0.2098214285714286
$ws.Range("R6").Value2 = 0.04910714285714286
$ws.Range("S6").Value2 = 0.3928571428571428
$ws.Range("B7").Value2 = 0.06572769953051644
$ws.Range("D7").Value2 = 0.004694835680751174
$ws.Range("F7").Value2 = 0.05164319248826291
$ws.Range("J7").Value2 = 0.136150234741784
$ws.Range("O7").Value2 = 0.02816901408450704
$ws.Range("Q7").Value2 = 0.1971830985915493
$ws.Range("R7").Value2 = 0.1173708920187793
$ws.Range("S7").Value2 = 0.3990610328638498
$ws.Range("B8").Value2 = 0.0625
$ws.Range("D8").Value2 = 0.02232142857142857
$ws.Range("F8").Value2 = 0.06919642857142858
$ws.Range("J8").Value2 = 0.09821428571428571
$ws.Range("O8").Value2 = 0.03348214285714286
$ws.Range("Q8").Value2 = 0.1741071428571428
$ws.Range("R8").Value2 = 0.08035714285714286
$ws.Range("S8").Value2 = 0.4598214285714285
$ws.Range("B9").Value2 = 0.0867579908675799
$ws.Range("D9").Value2 = 0.0136986301369863
$ws.Range("F9").Value2 = 0.0684931506849315
$ws.Range("J9").Value2 = 0.1187214611872146
$ws.Range("O9").Value2 = 0.0365296803652968
$ws.Range("Q9").Value2 = 0.1415525114155251
$ws.Range("R9").Value2 = 0.091324200913242
$ws.Range("S9").Value2 = 0.4429223744292237
$ws.Range("B10").Value2 = 0.1114348142753095
$ws.Range("D10").Value2 = 0.02549162418062637
$ws.Range("F10").Value2 = 0.06846321922796796
$ws.Range("J10").Value2 = 0.1223597960670066
$ws.Range("O10").Value2 = 0.01238164603058995
$ws.Range("Q10").Value2 = 0.2068463219227968
$ws.Range("R10").Value2 = 0.06627822286962855
$ws.Range("S10").Value2 = 0.3867443554260743
$ws.Range("G11").Value2 = 0.1523809523809524
$ws.Range("J11").Value2 = 0.0761904761904762
$ws.Range("K11").Value2 = 0.1873015873015873
$ws.Range("L11").Value2 = 0.5650793650793651
$ws.Range("S11").Value2 = 0.01904761904761905
$ws.Range("G12").Value2 = 0.7362637362637363
$ws.Range("J12").Value2 = 0.1703296703296703
$ws.Range("K12").Value2 = 0.01648351648351648
$ws.Range("L12").Value2 = 0.04945054945054945
$ws.Range("S12").Value2 = 0.02747252747252747
$ws.Range("G13").Value2 = 0.6888888888888889
$ws.Range("J13").Value2 = 0.2222222222222222
$ws.Range("S13").Value2 = 0.08888888888888889
$ws.Range("F15").Value2 = 0.007905138339920948
$ws.Range("H15").Value2 = 0.1185770750988142
$ws.Range("I15").Value2 = 0.08695652173913043
$ws.Range("J15").Value2 = 0.3636363636363636
$ws.Range("K15").Value2 = 0.05138339920948617
$ws.Range("M15").Value2 = 0.007905138339920948
$ws.Range("O15").Value2 = 0.09486166007905138
$ws.Range("S15").Value2 = 0.2687747035573123
$ws.Range("F16").Value2 = 0.01704545454545454
$ws.Range("H16").Value2 = 0.1363636363636364
$ws.Range("I16").Value2 = 0.07386363636363637
$ws.Range("J16").Value2 = 0.4829545454545455
$ws.Range("K16").Value2 = 0.1022727272727273
$ws.Range("M16").Value2 = 0.01136363636363636
$ws.Range("O16").Value2 = 0.04545454545454546
$ws.Range("S16").Value2 = 0.1306818181818182
$ws.Range("F17").Value2 = 0.01232032854209446
$ws.Range("H17").Value2 = 0.1581108829568789
$ws.Range("I17").Value2 = 0.06160164271047228
$ws.Range("J17").Value2 = 0.4229979466119096
$ws.Range("K17").Value2 = 0.09650924024640657
$ws.Range("M17").Value2 = 0.02464065708418891
$ws.Range("O17").Value2 = 0.06776180698151951
$ws.Range("S17").Value2 = 0.1560574948665298
$ws.Range("F18").Value2 = 0.01639344262295082
$ws.Range("H18").Value2 = 0.1311475409836066
$ws.Range("I18").Value2 = 0.1092896174863388
$ws.Range("J18").Value2 = 0.3770491803278688
$ws.Range("K18").Value2 = 0.08196721311475409
$ws.Range("M18").Value2 = 0.03278688524590164
$ws.Range("O18").Value2 = 0.08743169398907104
$ws.Range("S18").Value2 = 0.1639344262295082
$ws.Range("F19").Value2 = 0.01554054054054054
$ws.Range("H19").Value2 = 0.2033783783783784
$ws.Range("I19").Value2 = 0.0918918918918919
$ws.Range("J19").Value2 = 0.3648648648648649
$ws.Range("K19").Value2 = 0.1067567567567568
$ws.Range("M19").Value2 = 0.01486486486486487
$ws.Range("N19").Value2 = 0.003378378378378379
$ws.Range("O19").Value2 = 0.06148648648648649
$ws.Range("S19").Value2 = 0.1378378378378378
